# Auto commit at 2025-10-19 8:11:50.27
# Refresh the Metrics feed values and recalc the dependent "today" sheet.

$wb = $excel.ActiveWorkbook

# --- Update raw metric values on the "Metrics" sheet ---
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value = 246596.69
$metrics.Range("B3").Value = 202353.80000000002
$metrics.Range("B4").Value = 78538.899999999994
$metrics.Range("B5").Value = 9848
$metrics.Range("B6").Value = 4613728.16
$metrics.Range("B7").Value = 3892172.4699999997
$metrics.Range("B8").Value = 1349141.04
$metrics.Range("B9").Value = 178849
$metrics.Range("B10").Value = 33079051.960999824
$metrics.Range("B11").Value = 31167393.990000002
$metrics.Range("B12").Value = 11630849.930000002
$metrics.Range("B13").Value = 1276476

# Move the saved selection on "Metrics" to F15, matching the author's cursor
[void]$metrics.Range("F15").Select()

# --- Move the saved selection on "today" to F7 ---
$today = $wb.Worksheets.Item("today")
$today.Activate()
[void]$today.Range("F7").Select()

# Recalculate so the dependent formulas on "today" (and TODAY()-1) pick up
# the refreshed values.
$excel.Calculate()
